# GUI SKETCH V3 upload
# Rebuild the "구현" (implementation progress) sheet with the new rows/columns,
# add a new wide column F, set a page setup on it, and move the active
# selection over to the "테스트시나리오" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Start from a clean slate for this sheet's data - the rows below row 14 are
# being substantially re-laid-out (new rows inserted, new columns of data
# added), so clear out the old body and re-write every cell at its final
# target location.
$ws.Cells.ClearContents()

# --- BACKEND section (unchanged positions, now with a "구현률" % column) ---
$ws.Range("B3").Value = "BACKEND"

$ws.Range("E4").Value = "구현률"
$ws.Range("F4").Value = "담당자"

$ws.Range("B5").Value = "시나리오"
$ws.Range("C5").Value = "기능"
$ws.Range("D5").Value = "API"

$ws.Range("C6").Value = "회원가입/로그인"
$ws.Range("D6").Value = "로그인"
$ws.Range("E6").Value = 100

$ws.Range("D7").Value = "회원가입"
$ws.Range("E7").Value = 100

$ws.Range("D8").Value = "이메일 중복 확인"
$ws.Range("E8").Value = 100

$ws.Range("D9").Value = "회원 선택지역 저장"
$ws.Range("E9").Value = 100

$ws.Range("D10").Value = "회원 선택 지역 조회"
$ws.Range("E10").Value = 100

$ws.Range("D11").Value = "회원 선호 정보 저장"
$ws.Range("E11").Value = 100

$ws.Range("D12").Value = "회원 선호 정보 조회"
$ws.Range("E12").Value = 100

$ws.Range("D13").Value = "회원 퍼스널리티 저장"
$ws.Range("E13").Value = 100

$ws.Range("D14").Value = "회원 퍼스널리티 조회"
$ws.Range("E14").Value = 100

# --- "구현 Element" section header ---
$ws.Range("D22").Value = "구현 Element"

# --- FRONT section ---
$ws.Range("B32").Value = "FRONT"

$ws.Range("C33").Value = "로그인화면"
$ws.Range("E33").Value = 100

$ws.Range("D34").Value = "로그인 성공"
$ws.Range("D35").Value = "로그인 실패"

$ws.Range("C37").Value = "회원가입화면"
$ws.Range("E37").Value = 100

$ws.Range("C38").Value = "지역추천화면"

$ws.Range("C40").Value = "로딩화면"

$ws.Range("C42").Value = "메인화면"

$ws.Range("D43").Value = "아이템"
$ws.Range("E43").Value = 100

$ws.Range("D44").Value = "좋아요버튼"
$ws.Range("D45").Value = "싫어요 버튼"
$ws.Range("D46").Value = "새로고침버튼"
$ws.Range("E46").Value = 100
$ws.Range("D47").Value = "설문조사 화면"
$ws.Range("E44").Value = "?"
$ws.Range("E45").Value = "?"

$ws.Range("C49").Value = "상세페이지"

# --- BACKEND "지역추천" / "피드백" additions higher up in the sheet ---
$ws.Range("C15").Value = "지역추천"
$ws.Range("C19").Value = "피드백 기능"
$ws.Range("D20").Value = "피드백 후 새로고침 반영"

# --- DATA section ---
$ws.Range("B55").Value = "DATA"
$ws.Range("C56").Value = "Data input"
$ws.Range("C58").Value = "data output"
$ws.Range("D56").Value = "Json을 사용한 데이터 Input"
$ws.Range("D58").Value = "Json을 이용한 데이터 Output"

# New wide column for descriptions
$ws.Range("F1").ColumnWidth = 49.57

# Page setup for this sheet (portrait, letter/A4-class paper 9)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("D16").Select()

# --- Active sheet moves to "테스트시나리오" with a selected cell ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("C7").Select()
